$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.558.97"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "'3.837.75"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'523.59"
$ws.Range("E5").Value = "  +7.42%  "
$ws.Range("D6").Value = "'142.57"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "'0.602"
$ws.Range("E7").Value = "  -3.17%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.707"
$ws.Range("E9").Value = "  -4.84%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "  -7.07%  "
$ws.Range("D11").Value = "'0.0000324"
$ws.Range("E11").Value = "  -8.53%  "
$ws.Range("D12").Value = "'41.46"
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").Value = "'4.442.95"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").Value = "'10.08"
$ws.Range("E14").Value = "  -3.79%  "
$ws.Range("D15").Value = "'3.842.01"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "'20.90"
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'1.21"
$ws.Range("E18").Value = "  +5.69%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.134"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").Value = "'68.449.55"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'415.84"
$ws.Range("E21").Value = "  -4.16%  "
$ws.Range("D22").Value = "'3.38"
$ws.Range("E22").Value = "  -4.29%  "
$ws.Range("D23").Value = "'13.90"
$ws.Range("E23").Value = "  -5.24%  "
$ws.Range("D24").Value = "'86.20"
$ws.Range("E24").Value = "  -4.98%  "
$ws.Range("E25").Value = "  +5.38%  "
$ws.Range("D26").Value = "'11.31"
$ws.Range("E26").Value = "  -7.89%  "
$ws.Range("D27").Value = "'10.46"
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("D28").Value = "'35.82"
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("D29").Value = "'681.33"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("D30").Value = "'13.00"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "'0.124"
$ws.Range("E31").Value = "  -4.78%  "
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").Value = "'66.84"
$ws.Range("E33").Value = "  +9.26%  "
$ws.Range("E34").Value = "  +6.92%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "'0.0₃0844"
$ws.Range("E35").Value = "  -5.74%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.83"
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").Value = "'39.35"
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Value = "'3.18"
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").Value = "'2.72"
$ws.Range("E44").Value = "  -7.56%  "
$ws.Range("D45").Value = "'3.39"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.93"
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "'0.000270"
$ws.Range("E48").Value = "  +9.80%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'2.726.86"
$ws.Range("E49").Value = "  +12.56%  "
$ws.Range("D50").Value = "'143.69"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("E51").Value = "  -3.82%  "
